$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "Status" column (D) values for the first four bookings.
$ws.Range("D1").Value = "CHECK-OUT"
$ws.Range("D2").Value = "CHECK-OUT"
$ws.Range("D4").Value = "CHECK-IN"

# Append two new booking records as rows 5 and 6. Booking ID / Room Number
# are numeric-looking text in this sheet (like the existing rows), so write
# them through TEXT()+paste-values to keep them stored as text rather than
# being auto-converted to numbers.
$ws.Range("A5").Formula = '=TEXT(1545,"0")'
$ws.Range("A5").Copy()
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("B5").Value = "Raweeroj   Thongdee"
$ws.Range("C5").Formula = '=TEXT(1003,"0")'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)
$ws.Range("D5").Value = "Waiting"
$ws.Range("E5").Value = "24-04-2020 03:32:19"

$ws.Range("A6").Formula = '=TEXT(1573,"0")'
$ws.Range("A6").Copy()
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("B6").Value = "Raweeroj   Thongdee"
$ws.Range("C6").Formula = '=TEXT(2003,"0")'
$ws.Range("C6").Copy()
$ws.Range("C6").PasteSpecial(-4163)
$ws.Range("D6").Value = "Waiting"
$ws.Range("E6").Value = "24-04-2020 03:32:19"

# Widen column D slightly to fit the new content, matching the target layout.
$ws.Columns.Item(4).ColumnWidth = 10.8
